# Fruta / hortaliza, semanal
# Insert 3 new price-report rows (Palta Hass, "Especial"/"Primera"/"Segunda",
# origin "Perú", sold by the 10kg tray) right before the current row 811,
# pushing the existing rows 811:871 down to 814:874.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("811:813").Insert()

# Shared attributes for this product block (identical on every Palta/Hass row).
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$variedad  = "Hass"

$newRows = @(
    @{ Row = 811; Fecha = 45013; Calidad = "Especial"; Volumen = 50;  PMin = 31000; PMax = 31000; PProm = 31000; Unidad = "`$/bandeja 10 kilos"; Origen = "Perú"; PrecioKg = 3100; KgUnidad = 10 },
    @{ Row = 812; Fecha = 45013; Calidad = "Primera";  Volumen = 100; PMin = 29000; PMax = 29000; PProm = 29000; Unidad = "`$/bandeja 10 kilos"; Origen = "Perú"; PrecioKg = 2900; KgUnidad = 10 },
    @{ Row = 813; Fecha = 45013; Calidad = "Segunda";  Volumen = 100; PMin = 26000; PMax = 26000; PProm = 26000; Unidad = "`$/bandeja 10 kilos"; Origen = "Perú"; PrecioKg = 2600; KgUnidad = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $r.Fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.PMin
    $ws.Cells.Item($row, 15).Value2 = $r.PMax
    $ws.Cells.Item($row, 16).Value2 = $r.PProm
    $ws.Cells.Item($row, 17).Value2 = $r.Unidad
    $ws.Cells.Item($row, 18).Value2 = $r.Origen
    $ws.Cells.Item($row, 19).Value2 = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value2 = $r.KgUnidad
}
